$d = $word.ActiveDocument

# -------------------------------------------------------------------
# Edit 1: merge the "Strengths ..." two paragraphs into one paragraph
# with new wording, and set spacing-before on the paragraph.
# -------------------------------------------------------------------
$p1 = $null
$p2 = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "^Strengths including analytical") { $p1 = $p }
    if ($t -match "^diverse groups, makes me a valuable") { $p2 = $p }
}

if ($p1 -ne $null -and $p2 -ne $null) {
    $r1 = $d.Range($p1.Range.Start, $p2.Range.End)

    $xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/dummy.xml" pkg:contentType="application/xml">
    <pkg:xmlData>
      <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:pPr>
          <w:pStyle w:val="Heading1"/>
          <w:spacing w:before="165"/>
          <w:ind w:left="0"/>
          <w:rPr>
            <w:b w:val="0"/>
            <w:bCs w:val="0"/>
            <w:color w:val="000000"/>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b w:val="0"/>
            <w:bCs w:val="0"/>
            <w:color w:val="000000"/>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t xml:space="preserve">Strengths </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:b w:val="0"/>
            <w:bCs w:val="0"/>
            <w:color w:val="000000"/>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t>include</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:b w:val="0"/>
            <w:bCs w:val="0"/>
            <w:color w:val="000000"/>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t xml:space="preserve"> stringing sentences together from alphabet soup, </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:b w:val="0"/>
            <w:bCs w:val="0"/>
            <w:color w:val="000000"/>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t xml:space="preserve">a </w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:b w:val="0"/>
            <w:bCs w:val="0"/>
            <w:color w:val="000000"/>
            <w:sz w:val="20"/>
            <w:szCs w:val="20"/>
          </w:rPr>
          <w:t>seasoned "cat herder" and experience collaborating across diverse, distributed global teams.</w:t>
        </w:r>
      </w:p>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

    $r1.InsertXML($xml1)
} else {
    Write-Host "WARNING: could not locate Strengths paragraphs for edit 1"
}

# -------------------------------------------------------------------
# Edit 2: split the "Resolved a multitude ... Manager" paragraph
# into three bulleted paragraphs.
# -------------------------------------------------------------------
$p3 = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "^Resolved a multitude of issues with customer") { $p3 = $p }
}

if ($p3 -ne $null) {
    $r2 = $d.Range($p3.Range.Start, $p3.Range.End)

    $xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/dummy.xml" pkg:contentType="application/xml">
    <pkg:xmlData>
      <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00263422" w:rsidRPr="00F81DD6">
        <w:pPr>
          <w:pStyle w:val="BodyText"/>
          <w:numPr>
            <w:ilvl w:val="0"/>
            <w:numId w:val="7"/>
          </w:numPr>
          <w:spacing w:line="292" w:lineRule="auto"/>
          <w:ind w:right="4588"/>
        </w:pPr>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
          </w:rPr>
          <w:t>Resolved</w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
            <w:spacing w:val="-6"/>
          </w:rPr>
          <w:t xml:space="preserve"> </w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
          </w:rPr>
          <w:t>a</w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
            <w:spacing w:val="-6"/>
          </w:rPr>
          <w:t xml:space="preserve"> </w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
          </w:rPr>
          <w:t>multitude</w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
            <w:spacing w:val="-6"/>
          </w:rPr>
          <w:t xml:space="preserve"> </w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
          </w:rPr>
          <w:t>of</w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
            <w:spacing w:val="-6"/>
          </w:rPr>
          <w:t xml:space="preserve"> </w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
          </w:rPr>
          <w:t>issues</w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
            <w:spacing w:val="-6"/>
          </w:rPr>
          <w:t xml:space="preserve"> </w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
          </w:rPr>
          <w:t>with</w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
            <w:spacing w:val="-6"/>
          </w:rPr>
          <w:t xml:space="preserve"> </w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
          </w:rPr>
          <w:t>customer</w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
            <w:spacing w:val="-6"/>
          </w:rPr>
          <w:t xml:space="preserve"> </w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
          </w:rPr>
          <w:t>devices</w:t>
        </w:r>
      </w:p>
      <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00263422" w:rsidRPr="00F81DD6">
        <w:pPr>
          <w:pStyle w:val="BodyText"/>
          <w:numPr>
            <w:ilvl w:val="0"/>
            <w:numId w:val="7"/>
          </w:numPr>
          <w:spacing w:line="292" w:lineRule="auto"/>
          <w:ind w:right="4588"/>
        </w:pPr>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
          </w:rPr>
          <w:t xml:space="preserve">Received Sprint Service Technician </w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:t>Certification</w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
          </w:rPr>
          <w:t xml:space="preserve"> </w:t>
        </w:r>
      </w:p>
      <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00263422" w:rsidRPr="00F81DD6">
        <w:pPr>
          <w:pStyle w:val="BodyText"/>
          <w:numPr>
            <w:ilvl w:val="0"/>
            <w:numId w:val="7"/>
          </w:numPr>
          <w:spacing w:line="292" w:lineRule="auto"/>
          <w:ind w:right="4588"/>
        </w:pPr>
        <w:r w:rsidRPr="00F81DD6">
          <w:t>P</w:t>
        </w:r>
        <w:r w:rsidRPr="00F81DD6">
          <w:rPr>
            <w:rFonts w:hint="cs"/>
          </w:rPr>
          <w:t>romoted to Service Technician Manager</w:t>
        </w:r>
      </w:p>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

    $r2.InsertXML($xml2)
} else {
    Write-Host "WARNING: could not locate the Resolved-a-multitude paragraph for edit 2"
}


# -------------------------------------------------------------------
# Edit 3: "Default Paragraph Font" style should no longer be
# semi-hidden (best-effort; silently ignored if the host doesn't
# expose a writer for this particular flag).
# -------------------------------------------------------------------
try {
    $dpf = $d.Styles("Default Paragraph Font")
    $dpf.Hidden = $false
} catch {
}

Write-Host "Done. Paragraph count now: $($d.Paragraphs.Count)"
